$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Satisfied By" (column D) for a few requirements
$ws.Range("D5").Value = "Waste Transfer Subsystem"
$ws.Range("D7").Value = "Navigation Subsystem"
$ws.Range("D9").Value = "Waste Transfer Subsystem"
$ws.Range("D11").Value = "Waste Transfer Subsystem"
